# Cambiar formato de las fechas
# Convert the date values in column A of "Hoja2" from Excel date serials
# to plain text strings formatted as "YYYY-MM" (e.g. 2020-06, 2020-07, ...)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja2")

# Ordered list of new text values for rows 2..14 (matches the dates
# previously stored as the serials 43983, 44013, 44044, 44075, 44105,
# 44136, 44166, 44197, 44228, 44256, 44287, 44317, 44348 -- the first day
# of each month from June 2020 through June 2021).
$dates = @(
    "2020-06",
    "2020-07",
    "2020-08",
    "2020-09",
    "2020-10",
    "2020-11",
    "2020-12",
    "2021-01",
    "2021-02",
    "2021-03",
    "2021-04",
    "2021-05",
    "2021-06"
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 1)
    # Simply assign the text value; the cell keeps its existing style/number
    # format (a date format), but since the underlying value is now text
    # Excel stores it as a shared string rather than a date serial number.
    $cell.Value = $dates[$i]
}

# Update the active selection on Hoja2 to match the authored workbook state.
$ws.Activate()
$ws.Range("B11").Select()

$wb.Save()
